$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1734425452534, "adesh", "adesh@1.com", '$2a$10$qV6Jcj5LqwPcdZAr/Y4poOANzsfwtdL2.aev0IKJvnK0hVxStWKq.'),
    @(1734425777014, "gii",   "dd@1ss.com",  '$2a$10$HYsiBR56bxBVbrs8TeHxSeRzbCgC381Xc609pxDx56gw6wgLSkq/6'),
    @(1734426135137, "adesh", "adesh#@s.com",'$2a$10$xjYtaIQv/ipHVF.wEtu6QudHdqeOuqiCri9Hr0PAZ9gdysnOk2sUS'),
    @(1734426247575, "affa",  "adesh@q.com", '$2a$10$pOnD7iZI7apARgGmmiBIpOno0Dqo80oDEbIPE2ynWsEhOpq5WWgvi')
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
